$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (rows 2-51).
# Column D (Price) values are forced to Text format before assignment and then
# have formatting cleared, so numeric-looking strings (e.g. "6.00", "0.0000333")
# are preserved exactly as text instead of being auto-converted to numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.594.51"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.77%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.922.95"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.14%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +9.36%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.56"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.00%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.615"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.22%  "

# Row 8
$ws.Range("E8").Value = "  +0.11%  "

# Row 9
$ws.Range("E9").Value = "  -0.80%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.61%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000333"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.65%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.44"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.69%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.550.89"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.27%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.29"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.928.41"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.33%  "

# Row 16
$ws.Range("E16").Value = "  +8.41%  "

# Row 17
$ws.Range("E17").Value = "  -0.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.89"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.12%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.88"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.505.98"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.60%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "436.64"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.85%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.35"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.71%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.41"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.47%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +12.99%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.03"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.26%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.55"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.56%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.72"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.79%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.40"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.05%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "697.08"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.36%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.21"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.30%  "

# Row 31
$ws.Range("E31").Value = "  -2.68%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.83"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.55%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "68.80"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +12.74%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.450"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +14.10%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.22%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.36"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.41%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0839"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.42%  "

# Row 39
$ws.Range("E39").Value = "  -0.01%  "

# Row 40
$ws.Range("E40").Value = "  +0.01%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0486"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.51%  "

# Row 42
$ws.Range("E42").Value = "  +3.50%  "

# Row 43
$ws.Range("E43").Value = "  -9.17%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.95"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.13%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.13"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +10.87%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.142"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.48%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.36"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.33%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.99%  "

# Row 49
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "144.76"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.29%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0340"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.83%  "

# Row 51
$ws.Range("E51").Value = "  -3.26%  "
